$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "current balance" (الرصيد الحالي) column for all three item rows now
# reads "0:0" instead of the previous shortage figures ("-1:0" / "-23:0").
$ws.Range("H7").Value = "0:0"
$ws.Range("H8").Value = "0:0"
$ws.Range("H9").Value = "0:0"
